$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B3 gets a date value (2012-12-13) formatted the same as the other date
# cells in column B (copy the existing date format from B2 so it reuses
# the same style rather than minting a new number format).
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("B3").Value = 41256

# C4 gets a plain numeric value
$ws.Range("C4").Value = 3

# Update the sheet's active selection to H4
$ws.Range("H4").Select()
